$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-column style (s="2") from A328 down through the new date cells
$ws.Range("A328").Copy($ws.Range("A329:A343"))

$ws.Range("A329").Value = 44403
$ws.Range("B329").Value = 1
$ws.Range("C329").Value = 2
$ws.Range("D329").Value = 23.5654530458348

$ws.Range("A330").Value = 44404
$ws.Range("B330").Value = 0
$ws.Range("C330").Value = 2
$ws.Range("D330").Value = 23.5654530458348

$ws.Range("A331").Value = 44405
$ws.Range("B331").Value = 1
$ws.Range("C331").Value = 3
$ws.Range("D331").Value = 35.34817956875221

$ws.Range("A332").Value = 44406
$ws.Range("B332").Value = 0
$ws.Range("C332").Value = 3
$ws.Range("D332").Value = 35.34817956875221

$ws.Range("A333").Value = 44407
$ws.Range("B333").Value = 0
$ws.Range("C333").Value = 2
$ws.Range("D333").Value = 23.5654530458348

$ws.Range("A334").Value = 44408
$ws.Range("B334").Value = 0
$ws.Range("C334").Value = 2
$ws.Range("D334").Value = 23.5654530458348

$ws.Range("A335").Value = 44409
$ws.Range("B335").Value = 0
$ws.Range("C335").Value = 2
$ws.Range("D335").Value = 23.5654530458348

$ws.Range("A336").Value = 44410
$ws.Range("B336").Value = 2
$ws.Range("C336").Value = 3
$ws.Range("D336").Value = 35.34817956875221

$ws.Range("A337").Value = 44411
$ws.Range("B337").Value = 0
$ws.Range("C337").Value = 3
$ws.Range("D337").Value = 35.34817956875221

$ws.Range("A338").Value = 44412
$ws.Range("B338").Value = 0
$ws.Range("C338").Value = 2
$ws.Range("D338").Value = 23.5654530458348

$ws.Range("A339").Value = 44413
$ws.Range("B339").Value = 1
$ws.Range("C339").Value = 3
$ws.Range("D339").Value = 35.34817956875221

$ws.Range("A340").Value = 44414
$ws.Range("B340").Value = 0
$ws.Range("C340").Value = 3
$ws.Range("D340").Value = 35.34817956875221

$ws.Range("A341").Value = 44415
$ws.Range("B341").Value = 0
$ws.Range("C341").Value = 3
$ws.Range("D341").Value = 35.34817956875221

$ws.Range("A342").Value = 44416
$ws.Range("B342").Value = 2
$ws.Range("C342").Value = 5
$ws.Range("D342").Value = 58.91363261458702

$ws.Range("A343").Value = 44417
$ws.Range("B343").Value = 1
$ws.Range("C343").Value = 4
$ws.Range("D343").Value = 47.13090609166961
